$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source rows 2-5 (account numbers are text with leading zeros, so the
# cells must stay text-typed -- force text format first so the numeric-
# looking strings aren't auto-converted to numbers and lose their leading
# zeros) get reordered: LEILA moves from the top of the block to the
# bottom, and her balance changes from 58752.93 to 4000.
#
# Before:                      After:
#  2: 004208447 LEILA   58752.93     2: 004368468 AHMAD   21621.63
#  3: 004368468 AHMAD   21621.63     3: 005040864 ANDRE   7000
#  4: 005040864 ANDRE   7000         4: 004313254 GUSTAVO 4292
#  5: 004313254 GUSTAVO 4292         5: 004208447 LEILA   4000

$ws.Range("A2:A5").NumberFormat = "@"

$ws.Range("A2").Value = "004368468"
$ws.Range("B2").Value = "AHMAD"
$ws.Range("C2").Value = 21621.63

$ws.Range("A3").Value = "005040864"
$ws.Range("B3").Value = "ANDRE"
$ws.Range("C3").Value = 7000

$ws.Range("A4").Value = "004313254"
$ws.Range("B4").Value = "GUSTAVO"
$ws.Range("C4").Value = 4292

$ws.Range("A5").Value = "004208447"
$ws.Range("B5").Value = "LEILA"
$ws.Range("C5").Value = 4000
